$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the second row (the old "(m3/s) (MW) (MW) (GWh) (GWh) (GWh)" unit-label row
# at F2:K2), which shifts all the power-plant data rows below it up by one. Excel
# also drops the now-unused "(MW)"/"(GWh)" entries from the shared-string table and
# the now-empty trailing row (old row 117) falls off the used range automatically.
$ws.Rows.Item(2).Delete()

# The new top-left selection sits on the first data row (A2:K2), matching the
# post-edit view state instead of the old header-row selection (A1:K1).
$ws.Range("A2:K2").Select() | Out-Null
